$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 200, pushing existing rows 200..235 down to 201..236
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with the new data record
$ws.Cells.Item(200, 1).Value = 3
$ws.Cells.Item(200, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(200, 3).Value = "Coquimbo"
$ws.Cells.Item(200, 4).Value = 44505
$ws.Cells.Item(200, 5).Value = 5
$ws.Cells.Item(200, 6).Value = 100112040
$ws.Cells.Item(200, 7).Value = "Cilantro"
$ws.Cells.Item(200, 8).Value = "Sin especificar"
$ws.Cells.Item(200, 9).Value = "Primera"
$ws.Cells.Item(200, 10).Value = 130
$ws.Cells.Item(200, 11).Value = 2000
$ws.Cells.Item(200, 12).Value = 2000
$ws.Cells.Item(200, 13).Value = 2000
$ws.Cells.Item(200, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(200, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(200, 16).Value = 667
$ws.Cells.Item(200, 17).Value = 3
$ws.Cells.Item(200, 18).Value = "Hortaliza"
